# Archie's Setup Checklist - apply authored edits
$wb = $excel.ActiveWorkbook

$wsChecklist = $wb.Worksheets.Item("Setup Checklist")
$wsHarmonics = $wb.Worksheets.Item("Harmonics Calculator")

# --- Setup Checklist sheet edits ---

# Rows 20/21 swapped their task descriptions
$wsChecklist.Range("C20").Value = "Go to mixers override check/set directions"
$wsChecklist.Range("C21").Value = "Check Swash level and full positvie and negative"

# Row 22 wording tweak ("Set calibration" -> "Set Initial calibration")
$wsChecklist.Range("C22").Value = "Set Initial calibration % to 40 for Cyclic, Collective, Tail"

# Insert a new row 51 with an edgetx-x10-scripts link (pushes the Discord rows down one)
$wsChecklist.Rows.Item(51).Insert()
$wsChecklist.Range("C51").Value = "https://github.com/offer-shmuely/edgetx-x10-scripts"

# --- Harmonics Calculator sheet edits (input cells; downstream formulas auto-recalc) ---
$wsHarmonics.Range("D6").Value = 1494
$wsHarmonics.Range("D7").Value = 16
$wsHarmonics.Range("E7").Value = 151
$wsHarmonics.Range("D8").Value = 20
$wsHarmonics.Range("D21").Value = 24.9

# --- View state: active sheet moves from Setup Checklist to Harmonics Calculator ---
[void]$wsChecklist.Range("D48").Select()
[void]$wsHarmonics.Activate()
[void]$wsHarmonics.Range("C14").Select()
